# Applies the reviewed-deck edits:
#   1. Slide 20: correct the bullet text under "count()".
#   2. Refresh the cached "datetimeFigureOut" placeholder text (5/8/2019 ->
#      10/31/2019) on the slide master, every slide layout, and the notes
#      master - mirroring a Header & Footer "Apply to All" date refresh.

$p = $ppt.ActivePresentation

# --- 1. Slide-level text fix -------------------------------------------------
$slide20 = $p.Slides.Item(20)
for ($i = 1; $i -le $slide20.Shapes.Count; $i++) {
    $shp = $slide20.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "Her bir elemanın sayısını gösteren bir Long döner.") {
            $shp.TextFrame.TextRange.Text = "RDD eleman sayısını gösteren bir Long döner."
        }
    }
}

# --- helper: fix a "dt" placeholder shape's cached text ---------------------
function Update-DatePlaceholder($shapes, [string]$newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $cur = $shp.TextFrame.TextRange.Text
            if ($cur -eq "5/8/2019") {
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

# --- 2. Slide master ----------------------------------------------------
Update-DatePlaceholder $p.SlideMaster.Shapes "10/31/2019"

# --- 3. Every slide layout off the master --------------------------------
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder $layout.Shapes "10/31/2019"
}

# --- 4. Notes master ------------------------------------------------------
Update-DatePlaceholder $p.NotesMaster.Shapes "10/31/2019"
